$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "特变电工"
$ws.Range("C2").Value = "利欧股份"
$ws.Range("B3").Value = "中国西电"
$ws.Range("C3").Value = "海格通信"
$ws.Range("A4").Value = "海格通信"
$ws.Range("B4").Value = "海格通信"
$ws.Range("C4").Value = "国晟科技"
$ws.Range("A5").Value = "山子高科"
$ws.Range("B5").Value = "盈方微"
$ws.Range("C5").Value = "特变电工"
$ws.Range("A6").Value = "金风科技"
$ws.Range("B6").Value = "金风科技"
$ws.Range("C6").Value = "锋龙股份"
$ws.Range("A7").Value = "保变电气"
$ws.Range("B7").Value = "岩山科技"
$ws.Range("C7").Value = "山子高科"
$ws.Range("A8").Value = "湖南白银"
$ws.Range("B8").Value = "山子高科"
$ws.Range("C8").Value = "中国西电"
$ws.Range("A9").Value = "三变科技"
$ws.Range("C9").Value = "湖南白银"
$ws.Range("A10").Value = "汉缆股份"
$ws.Range("B10").Value = "保变电气"
$ws.Range("C10").Value = "三维通信"
$ws.Range("A11").Value = "岩山科技"
$ws.Range("B11").Value = "湖南白银"
$ws.Range("C11").Value = "通富微电"
$ws.Range("B12").Value = "航天发展"
$ws.Range("C12").Value = "航天发展"
$ws.Range("A13").Value = "蓝色光标"
$ws.Range("B13").Value = "白银有色"
$ws.Range("C13").Value = "金风科技"
$ws.Range("A14").Value = "白银有色"
$ws.Range("B14").Value = "中国电建"
$ws.Range("C14").Value = "平潭发展"
$ws.Range("A15").Value = "锋龙股份"
$ws.Range("B15").Value = "通富微电"
$ws.Range("C15").Value = "航天电子"
$ws.Range("A16").Value = "浙文互联"
$ws.Range("B16").Value = "新联电子"
$ws.Range("C16").Value = "白银有色"
$ws.Range("A17").Value = "盈方微"
$ws.Range("B17").Value = "汉缆股份"
$ws.Range("C17").Value = "岩山科技"
$ws.Range("A18").Value = "航天电子"
$ws.Range("B18").Value = "航天电子"
$ws.Range("C18").Value = "值得买"
$ws.Range("A19").Value = "中国电建"
$ws.Range("B19").Value = "锋龙股份"
$ws.Range("C19").Value = "神剑股份"
$ws.Range("A20").Value = "通富微电"
$ws.Range("B20").Value = "三变科技"
$ws.Range("C20").Value = "九鼎新材"
$ws.Range("A21").Value = "三维通信"
$ws.Range("B21").Value = "浙文互联"
$ws.Range("C21").Value = "蓝色光标"
